$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.329.14'
$ws.Range("E2").Value = '  -3.74%  '
$ws.Range("D3").Value = '2.984.63'
$ws.Range("E3").Value = '  -3.16%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '535.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.48'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.70%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '2.979.31'
$ws.Range("E8").Value = '  -3.11%  '
$ws.Range("E9").Value = '  -0.20%  '
$ws.Range("E10").Value = '  -5.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.08'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.445'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.89%  '
$ws.Range("E13").Value = '  -2.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.62'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.18%  '
$ws.Range("D15").Value = '3.465.07'
$ws.Range("E15").Value = '  -3.16%  '
$ws.Range("E16").Value = '  -1.76%  '
$ws.Range("D17").Value = '61.394.14'
$ws.Range("E17").Value = '  -3.59%  '
$ws.Range("D18").Value = '2.983.78'
$ws.Range("E18").Value = '  -3.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.59'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '464.74'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.15'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.669'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.89'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.88'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.88%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("E27").Value = '  -1.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.67'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.22%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.88'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.14%  '
$ws.Range("B31").Value = 'Mantle'
$ws.Range("C31").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.15'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '25.46'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.66%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.43'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.14%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '54.87'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.27'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.85'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.63%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '455.64'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -9.06%  '
$ws.Range("D38").Value = '3.148.09'
$ws.Range("E38").Value = '  -3.59%  '
$ws.Range("E39").Value = '  -2.00%  '
$ws.Range("E40").Value = '  +1.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0381'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.58%  '
$ws.Range("E42").Value = '  -1.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.43'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -9.48%  '
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '26.08'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.242'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.99'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.108'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.62%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '117.14'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.21%  '
$ws.Range("E50").Value = '  -8.83%  '
$ws.Range("E51").Value = '  +6.37%  '
